$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -14.1707
$ws.Range("D4").Value = -7.765700000000002
$ws.Range("E4").Value = 12.1494

$ws.Range("D5").Value = -8.1465

$ws.Range("C7").Value = -11.98530000000001

$ws.Range("D8").Value = -8.347799999999998

$ws.Range("E9").Value = 14.09590000000001

$ws.Range("C16").Value = -11.9979
$ws.Range("D16").Value = -8.805200000000003

$ws.Range("E18").Value = 12.8908
